$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp note (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 19 de Abril de 2020 a las 04:22"

# --- Refresh India stats (row 20) ---
$ws.Cells.Item(20, 2).Value = 15723
$ws.Cells.Item(20, 3).Value = 1
$ws.Cells.Item(20, 5).Value = 12739

# --- Refresh Mexico stats (row 34) ---
$ws.Cells.Item(34, 4).Value = 2627
$ws.Cells.Item(34, 5).Value = 4220

# --- Refresh Nueva Zelanda stats (row 71) ---
$ws.Cells.Item(71, 2).Value = 1431
$ws.Cells.Item(71, 3).Value = 9
$ws.Cells.Item(71, 4).Value = 912
$ws.Cells.Item(71, 5).Value = 507
$ws.Cells.Item(71, 7).Value = 1
$ws.Cells.Item(71, 8).Value = 12

# --- Insert Guatemala with fresh stats right after Kenia (row 117), ---
# --- pushing Sri Lanka and Mayotte down a row, and dropping the old ---
# --- Guatemala row that used to sit after Mayotte. ---
$ws.Cells.Item(118, 1).Value = "Guatemala"
$ws.Cells.Item(118, 2).Value = 257
$ws.Cells.Item(118, 3).Value = 22
$ws.Cells.Item(118, 4).Value = 21
$ws.Cells.Item(118, 5).Value = 229
$ws.Cells.Item(118, 6).Value = 3
$ws.Cells.Item(118, 7).Value = 0
$ws.Cells.Item(118, 8).Value = 7

$ws.Cells.Item(119, 1).Value = "Sri Lanka"
$ws.Cells.Item(119, 2).Value = 254
$ws.Cells.Item(119, 3).Value = 0
$ws.Cells.Item(119, 4).Value = 86
$ws.Cells.Item(119, 5).Value = 161
$ws.Cells.Item(119, 6).Value = 1
$ws.Cells.Item(119, 7).Value = 0
$ws.Cells.Item(119, 8).Value = 7

$ws.Cells.Item(120, 1).Value = "Mayotte"
$ws.Cells.Item(120, 2).Value = 254
$ws.Cells.Item(120, 3).Value = 0
$ws.Cells.Item(120, 4).Value = 117
$ws.Cells.Item(120, 5).Value = 133
$ws.Cells.Item(120, 6).Value = 6
$ws.Cells.Item(120, 7).Value = 0
$ws.Cells.Item(120, 8).Value = 4

# --- Refresh Jamaica stats (row 126) ---
$ws.Cells.Item(126, 2).Value = 173
$ws.Cells.Item(126, 3).Value = 10
$ws.Cells.Item(126, 5).Value = 143
